# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort "Periodo Mora" column into ascending order (2304, 2305, 2306).
# Row 17 (2305) is already in place; only rows 16 and 18 swap values.
$ws.Range("E16").Value = "2304"
$ws.Range("E18").Value = "2306"

# Update "Salario Basico" (column G) for all three worker periods.
$ws.Range("G16").Value = 1778930
$ws.Range("G17").Value = 1778930
$ws.Range("G18").Value = 1778930
